# Implement WorksheetWithColumnHeader setCellValue and updateRow methods
#
# Adds a third worksheet, "ColumnHeaderWithNoData", to the workbook: a
# fixture that has the usual 4-column header row (Name / Age / Home
# Country / Occupation) but no data rows beneath it - used to exercise
# WorksheetWithColumnHeader.setCellValue()/updateRow() against a sheet
# that starts out empty.

$wb = $excel.ActiveWorkbook

# "ColumnHeader" already has the exact header row we need (and the column
# widths / page setup the new fixture should share), so duplicate it
# instead of building a sheet from scratch - this also makes the new
# sheet reuse the existing shared-string entries for the header labels
# rather than creating duplicates.
$source = $wb.Worksheets.Item("ColumnHeader")
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)

# Copy(Before, After) - place the duplicate after the last existing sheet
# so it lands at the end of the tab strip; Excel makes the new copy the
# active sheet, which is what we want here too.
$null = $source.Copy($null, $lastSheet)

$newSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$newSheet.Name = "ColumnHeaderWithNoData"

# Strip the copied data rows, keeping only the header row (row 1) so the
# sheet represents "header present, no data".
$usedRows = $newSheet.UsedRange.Rows.Count
if ($usedRows -gt 1) {
    $null = $newSheet.Rows("2:" + $usedRows).Delete()
}

$null = $newSheet.Range("B5").Select()
